# Insert a new data row right before the current row 263, shifting rows
# 263-344 down to 264-345 (Excel copies formatting from the row above on
# insert, which matches the existing date-format style in column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(263).Insert()

# Fill the newly inserted row 263 with the new record's data.
$ws.Cells.Item(263,1).Value  = 7
$ws.Cells.Item(263,2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(263,3).Value  = 'Ñuble'
$ws.Cells.Item(263,4).Value  = 45120
$ws.Cells.Item(263,5).Value  = 16
$ws.Cells.Item(263,6).Value  = 'Fruta'
$ws.Cells.Item(263,7).Value  = 100108
$ws.Cells.Item(263,8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(263,9).Value  = 100108005
$ws.Cells.Item(263,10).Value = 'Piña'
$ws.Cells.Item(263,11).Value = 'Caramelo'
$ws.Cells.Item(263,12).Value = 'Primera'
$ws.Cells.Item(263,13).Value = 80
$ws.Cells.Item(263,14).Value = 23000
$ws.Cells.Item(263,15).Value = 23000
$ws.Cells.Item(263,16).Value = 23000
$ws.Cells.Item(263,17).Value = '$/caja 12 unidades'
$ws.Cells.Item(263,18).Value = 'Ecuador'
$ws.Cells.Item(263,19).Value = 1917
$ws.Cells.Item(263,20).Value = 12
